$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U8").Value = 3
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 1
$ws.Range("AG8").Value = 1
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 3
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 5
$ws.Range("W9").Value = 6
$ws.Range("AA9").Value = 1
$ws.Range("AB9").Value = 2
$ws.Range("AC9").Value = 3
$ws.Range("AD9").Value = 4
$ws.Range("AE9").Value = 5
$ws.Range("AF9").Value = 6
$ws.Range("AG9").Value = 7
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 5
$ws.Range("Q10").Value = "frwd"
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 4
$ws.Range("V10").Value = 7
$ws.Range("W10").Value = 10
$ws.Range("Z10").Value = "frwd"
$ws.Range("AA10").Value = 1
$ws.Range("AB10").Value = 2
$ws.Range("AC10").Value = 5
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 11
$ws.Range("AF10").Value = 14
$ws.Range("AG10").Value = 17
$ws.Range("G11").Value = "frwd"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 6
$ws.Range("L11").Value = 9
$ws.Range("Q11").Value = "keep"
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 5
$ws.Range("T11").Value = 8
$ws.Range("U11").Value = 11
$ws.Range("Z11").Value = "keep"
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 3
$ws.Range("AC11").Value = 6
$ws.Range("AD11").Value = 9
$ws.Range("AE11").Value = 12
$ws.Range("AF11").Value = 15
$ws.Range("AG11").Value = 18
$ws.Range("G12").Value = "keep"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 7
$ws.Range("K12").Value = 10
$ws.Range("Q12").Value = "prev"
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 9
$ws.Range("T12").Value = 12
$ws.Range("Z12").Value = "prev"
$ws.Range("AA12").Value = 4
$ws.Range("AB12").Value = 7
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 13
$ws.Range("AE12").Value = 16
$ws.Range("AF12").Value = 19
$ws.Range("G13").Value = "prev"
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 8
$ws.Range("J13").Value = 11
$ws.Range("AA14").Value = 1
$ws.Range("Z15").Value = "frwd"
$ws.Range("AA15").Value = 1
$ws.Range("Z16").Value = "keep"
$ws.Range("AA16").Value = 0
$ws.Range("Z17").Value = "prev"

$ws.Rows.Item(8).RowHeight = 18
$ws.Rows.Item(9).RowHeight = 18
$ws.Rows.Item(10).RowHeight = 18
$ws.Rows.Item(11).RowHeight = 18
$ws.Rows.Item(12).RowHeight = 18
$ws.Rows.Item(13).RowHeight = 18
$ws.Rows.Item(14).RowHeight = 18
$ws.Rows.Item(15).RowHeight = 18
$ws.Rows.Item(16).RowHeight = 18
$ws.Rows.Item(17).RowHeight = 18

$ws.Range("AB8").Select()
